$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 436, shifting existing rows 436:529 down to 437:530.
$ws.Rows.Item(436).Insert()

# Populate the newly inserted row 436 with the new data record.
$ws.Cells.Item(436, 1).Value = 3
$ws.Cells.Item(436, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(436, 3).Value = "Coquimbo"
$ws.Cells.Item(436, 4).Value = 45209
$ws.Cells.Item(436, 5).Value = 5
$ws.Cells.Item(436, 6).Value = 100112001
$ws.Cells.Item(436, 7).Value = "Berenjena"
$ws.Cells.Item(436, 8).Value = "Sin especificar"
$ws.Cells.Item(436, 9).Value = "Primera"
$ws.Cells.Item(436, 10).Value = 65
$ws.Cells.Item(436, 11).Value = 8000
$ws.Cells.Item(436, 12).Value = 8000
$ws.Cells.Item(436, 13).Value = 8000
$ws.Cells.Item(436, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(436, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(436, 16).Value = 133
$ws.Cells.Item(436, 17).Value = 60
$ws.Cells.Item(436, 18).Value = "Hortaliza"
